$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.129642486572266
$ws.Range("B1").Value = 4.237920761108398
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.011576175689697
$ws.Range("E1").Value = 2.399283885955811
